$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: delete the text found by an exact, case-sensitive search and
# retype it, then reapply Bold / a character style to the freshly inserted
# span. Deleting + retyping (instead of Find&Replace / Range.Text=) is what
# lets the run coalesce the same way the rest of the (unedited) runs in the
# paragraph do, and also drops any <w:proofErr/> markers that wrapped the
# original text.
# ---------------------------------------------------------------------------
function Replace-Bold([string]$searchText, [string]$newText) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($searchText)
    if ($idx -lt 0) { Write-Output ("NOT FOUND: " + $searchText); return }
    $len = $searchText.Length
    $r = $d.Range($idx, $idx + $len)
    $r.Delete()
    $r2 = $d.Range($idx, $idx)
    $r2.InsertBefore($newText)
    $d.Range($idx, $idx + $newText.Length).Bold = 1
}

function Replace-Code([string]$searchText, [string]$newText) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($searchText)
    if ($idx -lt 0) { Write-Output ("NOT FOUND: " + $searchText); return }
    $len = $searchText.Length
    $r = $d.Range($idx, $idx + $len)
    $r.Delete()
    $r2 = $d.Range($idx, $idx)
    $r2.InsertBefore($newText)
    $d.Range($idx, $idx + $newText.Length).Style = "Code"
}

# ---------------------------------------------------------------------------
# 1) Title: "... wizards for Maya 2014 .NET API"
#         -> "... wizards for Maya 2015, 2016 and 2017 .NET API"
#    The "_GoBack" bookmark sits right after "wizards" and must stay put,
#    ending up between "2016 and " and "2017".
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")

# Text immediately following the bookmark is " for Maya 2014 .NET API".
# Turn " for Maya 2014" (everything up to, but not including, the trailing
# space before ".NET") into "2017" - this rewrites the run(s) right after
# the bookmark without ever touching the bookmark's own position.
$bmEnd = $bm.Range.End
$oldTail = " for Maya 2014"
$rTail = $d.Range($bmEnd, $bmEnd + $oldTail.Length)
if ($rTail.Text -ne $oldTail) {
    Write-Output ("unexpected tail text: [" + $rTail.Text + "]")
}
$rTail.Text = "2017"

# Insert the new lead-in text right before the bookmark (this does not
# disturb the bookmark, which stays anchored after the inserted text).
$bm2 = $d.Bookmarks("_GoBack")
$bmStart2 = $bm2.Range.Start
$rLead = $d.Range($bmStart2, $bmStart2)
$rLead.InsertBefore(" for Maya 2015, 2016 and ")

# ---------------------------------------------------------------------------
# 2) "Copy the Maya CSharp plug-in.zip to the ... Visual Studio
#     2010/Templates/ProjectTemplates/Visual C# directory."
#    -> bump 2010 to 2012 and drop the stray proofErr-induced run splits.
# ---------------------------------------------------------------------------
Replace-Bold "Maya CSharp plug-in." "Maya CSharp plug-in."
Replace-Code "2010/Templates/ProjectTemplates/Visual C#" "2012/Templates/ProjectTemplates/Visual C#"

# ---------------------------------------------------------------------------
# 3) "Copy the Maya CSharp Command with Undo.zip, Maya CSharp Command.zip,
#     and Maya CSharp Node.zip files to the ... Visual Studio
#     2010/Templates/ItemTemplates/Visual C# directory."
#    -> bump 2010 to 2012 and drop the stray proofErr-induced run splits.
# ---------------------------------------------------------------------------
Replace-Bold "Maya CSharp Command with Undo.zip" "Maya CSharp Command with Undo.zip"
Replace-Bold "Maya CSharp Command.zip" "Maya CSharp Command.zip"
Replace-Bold "Maya CSharp Node.zip" "Maya CSharp Node.zip"
Replace-Code "2010/Templates/ItemTemplates/Visual C#" "2012/Templates/ItemTemplates/Visual C#"

$d.Save()
